$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.946.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.90%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.25%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.018"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.59%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.017"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.47%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4690"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.71%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3909"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.56%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07962"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.007"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.84%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.63"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.23%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.902.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.35%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.931"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.080"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.41%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.020"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.55%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06775"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.44%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001042"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "

# Row 21
$ws.Range("E21").Value = "  +1.57%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.961.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.468"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.83%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.356"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.47%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.113.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.71%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.41%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.066"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.51%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.467"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.35%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.02"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.27%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09499"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.69%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9577"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.663"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.306"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.350"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.98%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06128"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02235"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.23%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.220"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.28%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.121"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5901"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1883"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.15"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.90%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5634"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.56%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.389"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.915"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.11%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06877"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.86%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.63"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.86%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.061"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.08%  "
